# Update gh-pages output data for 丽水-漫展信息 workbook.
# Changes apply identically to the "展览" and "全部类型" sheets:
#   F2: 168 -> 169   (想去人数 for row 2)
#   G2: 50  -> 60    (最低票价 for row 2)
#   G3: 29.9 -> 55   (最低票价 for row 3)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 169
    $ws.Range("G2").Value = 60
    $ws.Range("G3").Value = 55
}
